# #5: property boat&car done
# Fix the "汽車" (car) sheet: label the existing "capacity" column properly
# and append the standard trailing metadata columns (property_category,
# category, date, legislator_name, legislator_id, source_file, index) that
# already exist on the other property sheets (land/building/stock/etc).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("汽車")

# --- Header row (row 1) ---------------------------------------------------
# Row 1 had been mistakenly populated with the first data record's values
# instead of column headers; replace it with the standard header labels
# used by every other property sheet (land/building/stock/etc), and add
# the "capacity" header that was missing outright (column C used to show
# a bare number with no label at all).
$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "capacity"
$ws.Range("D1").Value = "owner"
$ws.Range("E1").Value = "register_date"
$ws.Range("F1").Value = "register_reason"
$ws.Range("G1").Value = "acquire_value"

$ws.Range("H1").Value = "property_category"
$ws.Range("I1").Value = "category"
$ws.Range("J1").Value = "date"
$ws.Range("K1").Value = "legislator_name"
$ws.Range("L1").Value = "legislator_id"
$ws.Range("M1").Value = "source_file"
$ws.Range("N1").Value = "index"

# Match the look (bold, centered, bordered header) of the rest of row 1.
$ws.Range("B1").Copy()
$ws.Range("H1:N1").PasteSpecial(-4122)

# --- Data row (row 2) ------------------------------------------------------
# Match the look (plain, unbordered) of the rest of row 2 first ...
$ws.Range("B2").Copy()
$ws.Range("H2:N2").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ... then fill in the values.
$ws.Range("H2").Value = "land"
$ws.Range("I2").Value = "normal"
# "2011-11-22" looks like a date to Excel's auto-detection, so force the
# cell to text first to keep it a literal string (matches every other
# sheet's "date" column, which stores this as plain text, not a date
# serial number).
$ws.Range("J2").NumberFormat = "@"
$ws.Range("J2").Value = "2011-11-22"
$ws.Range("K2").Value = "楊麗環"
$ws.Range("L2").Value = 960
$ws.Range("M2").Value = "tmpf3421"
$ws.Range("N2").Value = 30
